$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (DAMSLTag, DialogAct)
$updates = @{
    3  = @("sd", "Statement-non-opinion")
    8  = @("%", "Uninterpretable")
    9  = @("sd", "Statement-non-opinion")
    11 = @("sd", "Statement-non-opinion")
    15 = @("sd", "Statement-non-opinion")
    16 = @("ba", "Appreciation")
    25 = @("b", "Acknowledge (Backchannel)")
    28 = @("aa", "Agree/Accept")
    30 = @("aa", "Agree/Accept")
    33 = @("ba", "Appreciation")
    34 = @("aa", "Agree/Accept")
    43 = @("ba", "Appreciation")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
